$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update AgTests (H) / AgPosit (I) values for existing rows 282-301 ---
$ws.Cells.Item(282, 8).Value = 46903

$ws.Cells.Item(285, 8).Value = 40950
$ws.Cells.Item(285, 9).Value = 3441

$ws.Cells.Item(286, 8).Value = 55412
$ws.Cells.Item(286, 9).Value = 4260

$ws.Cells.Item(287, 8).Value = 57336
$ws.Cells.Item(287, 9).Value = 3911

$ws.Cells.Item(288, 8).Value = 57178
$ws.Cells.Item(288, 9).Value = 4008

$ws.Cells.Item(289, 8).Value = 64737
$ws.Cells.Item(289, 9).Value = 3666

$ws.Cells.Item(292, 8).Value = 80874
$ws.Cells.Item(292, 9).Value = 7159

$ws.Cells.Item(293, 8).Value = 82398
$ws.Cells.Item(293, 9).Value = 5860

$ws.Cells.Item(294, 8).Value = 90608
$ws.Cells.Item(294, 9).Value = 4941

$ws.Cells.Item(299, 8).Value = 62461
$ws.Cells.Item(299, 9).Value = 6487

$ws.Cells.Item(300, 8).Value = 68069
$ws.Cells.Item(300, 9).Value = 6628

$ws.Cells.Item(301, 8).Value = 64702
$ws.Cells.Item(301, 9).Value = 5164

# --- Append new daily rows 302-305 ---
$newRows = @(
    @(302, 44196, 184508, 127190, 55068, 16479, 4965, 2250, 63225, 4620),
    @(303, 44197, 186244, 128285, 55709, 4954, 1736, 2250, 9566, 678),
    @(304, 44198, 187463, 129994, 55152, 4288, 1219, 2317, 6388, 473),
    @(305, 44199, 188099, 130897, 54681, 3111, 636, 2521, 3404, 312)
)

foreach ($rowSpec in $newRows) {
    $r = $rowSpec[0]
    for ($c = 1; $c -le 9; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowSpec[$c]
    }
}
